$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text cells to stay as text (not auto-converted to
# numbers) by temporarily formatting as Text, then reset the style back to
# Normal so no stray number-format style is left attached to the cell.

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0.07351332"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.16822691"
$ws.Range("C2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "0.057512276"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0.14393607"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = 76
$ws.Range("E3").Value = 2

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0.18794297"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0.3455809"
$ws.Range("C4").Style = "Normal"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "0.1403265"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "0.28761235"
$ws.Range("C5").Style = "Normal"
